$d = $word.ActiveDocument

# Locate the "{% for citation in literature %}" paragraph and the
# "{{citation}}{% endfor %}" paragraph that follows it, then build a
# range spanning both paragraphs (including their paragraph marks).
$rngStart = $d.Content
$null = $rngStart.Find.Execute("{% for citation in literature %}")
$rngStart.Expand(4)  # wdParagraph - extend to include the paragraph mark

$rngEnd = $d.Content
$null = $rngEnd.Find.Execute("{{citation}}{% endfor %}")
$rngEnd.Expand(4)  # wdParagraph - extend to include the paragraph mark

$target = $d.Range($rngStart.Start, $rngEnd.End)

# Replace the two paragraphs with the six literature-reference paragraphs
# (styled "citation1") plus two trailing empty "citation1" paragraphs.
$xml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.integration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.absorption</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.</w:t></w:r><w:r><w:t>solution</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.refinement</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>literature.</w:t></w:r><w:r><w:t>ccdc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/></w:pPr><w:r><w:t>{{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>literature.</w:t></w:r><w:r><w:t>finalcif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="citation1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)

# The flat-OPC InsertXML path in this host drops <w:ind> set via pPr, so
# reapply the 360-twip (18pt) left indent on the two trailing empty
# paragraphs through the Paragraph.LeftIndent COM property instead.
$count = $d.Paragraphs.Count
$d.Paragraphs.Item($count - 1).LeftIndent = 18
$d.Paragraphs.Item($count).LeftIndent = 18
